$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5686
$ws1.Range("F5").Value = 313
$ws1.Range("F6").Value = 852
$ws1.Range("F7").Value = 64
$ws1.Range("F8").Value = 383
$ws1.Range("F11").Value = 22

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 22

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5686
$ws4.Range("F5").Value = 313
$ws4.Range("F6").Value = 852
$ws4.Range("F7").Value = 64
$ws4.Range("F9").Value = 383
$ws4.Range("F12").Value = 22
$ws4.Range("F13").Value = 22
